$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9, column A: tiny correction of the stored date/time serial value.
$ws.Range("A9").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A9").Value2 = 45873.70865909722

# New row 10 - same shape/format as the preceding rows.
$ws.Range("A10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A10").Value2 = 45873.75021946881
$ws.Range("B10").Value = 2025
$ws.Range("C10").Value = 15
$ws.Range("D10").Value = 18.56
$ws.Range("E10").Value = 78.36
$ws.Range("F10").Value = 14.06
$ws.Range("G10").Value = 7.74
$ws.Range("H10").Value = "ESE"
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = "18:00:18"
